# Auto-generated edit script: updates market-price-derived value cells
# across the 8 Leve profit worksheets, per the scheduled runner diff.
$wb = $excel.ActiveWorkbook

# ================= Sheet: ALC =================
$ws = $wb.Worksheets.Item("ALC")

# Row 121
$ws.Range("H121").Value = 5049.6665
$ws.Range("J121").Value = 5049.6665
$ws.Range("L121").Value = 15148.9995
$ws.Range("N121").Value = -18642.9995

# Row 132
$ws.Range("H132").Value = 6597.6665
$ws.Range("I132").Value = 1603.5454
$ws.Range("J132").Value = 20331.5
$ws.Range("K132").Value = 4810.6362
$ws.Range("L132").Value = 60994.5
$ws.Range("M132").Value = -2280.6362
$ws.Range("N132").Value = -66054.5

# ================= Sheet: ARM =================
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 2839.2727
$ws.Range("I2").Value = 2576.8572
$ws.Range("K2").Value = 2576.8572
$ws.Range("M2").Value = -2463.8572

# Row 40
$ws.Range("H40").Value = 27000
$ws.Range("J40").Value = 27000
$ws.Range("L40").Value = 27000
$ws.Range("N40").Value = -27352

# Row 45
$ws.Range("H45").Value = 14731.88
$ws.Range("I45").Value = 13824.647
$ws.Range("K45").Value = 13824.647
$ws.Range("M45").Value = -13447.647

# Row 52
$ws.Range("H52").Value = 60780
$ws.Range("J52").Value = 60780
$ws.Range("L52").Value = 60780
$ws.Range("N52").Value = -61416

# Row 61
$ws.Range("H61").Value = 3109.8
$ws.Range("I61").Value = 3109.8
$ws.Range("K61").Value = 3109.8
$ws.Range("M61").Value = -2897.8

# Row 74
$ws.Range("H74").Value = 1883.8572
$ws.Range("I74").Value = 1264.1111
$ws.Range("K74").Value = 1264.1111
$ws.Range("M74").Value = -390.1111000000001

# Row 77
$ws.Range("H77").Value = 1883.8572
$ws.Range("I77").Value = 1264.1111
$ws.Range("K77").Value = 6320.5555
$ws.Range("M77").Value = -1952.5555

# Row 116
$ws.Range("H116").Value = 2839.2727
$ws.Range("I116").Value = 2576.8572
$ws.Range("K116").Value = 2576.8572
$ws.Range("M116").Value = -282.8571999999999

# Row 132
$ws.Range("H132").Value = 2978.4
$ws.Range("I132").Value = 2054
$ws.Range("J132").Value = 6676
$ws.Range("K132").Value = 6162
$ws.Range("L132").Value = 20028
$ws.Range("M132").Value = -3632
$ws.Range("N132").Value = -25088

# Row 136
$ws.Range("H136").Value = 3109.8
$ws.Range("I136").Value = 3109.8
$ws.Range("K136").Value = 9329.400000000001
$ws.Range("M136").Value = -6779.400000000001

# ================= Sheet: BSM =================
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 2839.2727
$ws.Range("I3").Value = 2576.8572
$ws.Range("K3").Value = 2576.8572
$ws.Range("M3").Value = -2462.8572

# Row 19
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# Row 105
$ws.Range("H105").Value = 4167998
$ws.Range("I105").Value = 6251382
$ws.Range("J105").Value = 1229.75
$ws.Range("K105").Value = 6251382
$ws.Range("L105").Value = 1229.75
$ws.Range("M105").Value = -6249635
$ws.Range("N105").Value = -4723.75

# ================= Sheet: CRP =================
$ws = $wb.Worksheets.Item("CRP")

# Row 15
$ws.Range("H15").Value = 2000
$ws.Range("I15").Value = 2000
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 2000
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -1830
$ws.Range("N15").ClearContents()

# Row 31
$ws.Range("H31").Value = 2759.5
$ws.Range("I31").Value = 2759.5
$ws.Range("K31").Value = 2759.5
$ws.Range("M31").Value = -2464.5

# Row 34
$ws.Range("H34").Value = 2759.5
$ws.Range("I34").Value = 2759.5
$ws.Range("K34").Value = 2759.5
$ws.Range("M34").Value = -2557.5

# Row 38
$ws.Range("H38").Value = 8837.556
$ws.Range("I38").Value = 7038
$ws.Range("J38").Value = 9062.5
$ws.Range("K38").Value = 7038
$ws.Range("L38").Value = 9062.5
$ws.Range("M38").Value = -6661
$ws.Range("N38").Value = -9816.5

# Row 46
$ws.Range("H46").Value = 8837.556
$ws.Range("I46").Value = 7038
$ws.Range("J46").Value = 9062.5
$ws.Range("K46").Value = 7038
$ws.Range("L46").Value = 9062.5
$ws.Range("M46").Value = -6827
$ws.Range("N46").Value = -9484.5

# Row 58
$ws.Range("H58").Value = 6914.6
$ws.Range("I58").Value = 7791.6
$ws.Range("K58").Value = 7791.6
$ws.Range("M58").Value = -7588.6

# Row 132
$ws.Range("H132").Value = 2083.7334
$ws.Range("I132").Value = 1970.5834
$ws.Range("J132").Value = 2536.3333
$ws.Range("K132").Value = 5911.7502
$ws.Range("L132").Value = 7608.999899999999
$ws.Range("M132").Value = -3381.7502
$ws.Range("N132").Value = -12668.9999

# Row 136
$ws.Range("H136").Value = 6914.6
$ws.Range("I136").Value = 7791.6
$ws.Range("K136").Value = 23374.8
$ws.Range("M136").Value = -20824.8

# ================= Sheet: CUL =================
$ws = $wb.Worksheets.Item("CUL")

# Row 107
$ws.Range("H107").Value = 1236.8125
$ws.Range("I107").Value = 847.5
$ws.Range("K107").Value = 2542.5
$ws.Range("M107").Value = -622.5

# Row 131
$ws.Range("H131").Value = 8830.134
$ws.Range("J131").Value = 2519.8
$ws.Range("L131").Value = 7559.400000000001
$ws.Range("N131").Value = -17639.4

# ================= Sheet: GSM =================
$ws = $wb.Worksheets.Item("GSM")

# Row 2
$ws.Range("H2").Value = 385.52942
$ws.Range("I2").Value = 351.55554
$ws.Range("K2").Value = 351.55554
$ws.Range("M2").Value = -238.55554

# Row 21
$ws.Range("H21").Value = 1000
$ws.Range("J21").Value = 1000
$ws.Range("L21").Value = 1000
$ws.Range("N21").Value = -1346

# Row 30
$ws.Range("H30").Value = 1000
$ws.Range("J30").Value = 1000
$ws.Range("L30").Value = 1000
$ws.Range("N30").Value = -1210

# Row 32
$ws.Range("H32").Value = 83333.336
$ws.Range("J32").Value = 83333.336
$ws.Range("L32").Value = 83333.336
$ws.Range("N32").Value = -83925.336

# Row 33
$ws.Range("H33").Value = 19642.857
$ws.Range("J33").Value = 21500
$ws.Range("L33").Value = 21500
$ws.Range("N33").Value = -22004

# Row 107
$ws.Range("H107").Value = 4316.2144
$ws.Range("I107").Value = 4263.857
$ws.Range("J107").Value = 4368.5713
$ws.Range("K107").Value = 4263.857
$ws.Range("L107").Value = 4368.5713
$ws.Range("M107").Value = -2343.857
$ws.Range("N107").Value = -8208.5713

# Row 132
$ws.Range("H132").Value = 2038
$ws.Range("I132").Value = 2085.125
$ws.Range("K132").Value = 6255.375
$ws.Range("M132").Value = -3725.375

# ================= Sheet: LTW =================
$ws = $wb.Worksheets.Item("LTW")

# Row 31
$ws.Range("H31").Value = 2116
$ws.Range("I31").Value = 1814.3846
$ws.Range("J31").Value = 2900.2
$ws.Range("K31").Value = 1814.3846
$ws.Range("L31").Value = 2900.2
$ws.Range("M31").Value = -1566.3846
$ws.Range("N31").Value = -3396.2

# Row 33
$ws.Range("H33").Value = 28600
$ws.Range("I33").Value = 28600
$ws.Range("K33").Value = 28600
$ws.Range("M33").Value = -28310

# Row 40
$ws.Range("H40").Value = 5575
$ws.Range("I40").Value = 5691.091
$ws.Range("J40").Value = 5415.375
$ws.Range("K40").Value = 5691.091
$ws.Range("L40").Value = 5415.375
$ws.Range("M40").Value = -5555.091
$ws.Range("N40").Value = -5687.375

# Row 41
$ws.Range("H41").Value = 15000
$ws.Range("I41").Value = 10000
$ws.Range("J41").Value = 20000
$ws.Range("K41").Value = 10000
$ws.Range("L41").Value = 20000
$ws.Range("M41").Value = -9562
$ws.Range("N41").Value = -20876

# Row 55
$ws.Range("H55").Value = 318.3
$ws.Range("I55").Value = 225.42857
$ws.Range("K55").Value = 225.42857
$ws.Range("M55").Value = -52.42857000000001

# Row 132
$ws.Range("H132").Value = 7765.278
$ws.Range("I132").Value = 8572.333000000001
$ws.Range("K132").Value = 25716.999
$ws.Range("M132").Value = -23186.999

# ================= Sheet: WVR =================
$ws = $wb.Worksheets.Item("WVR")

# Row 17
$ws.Range("H17").Value = 4001
$ws.Range("I17").Value = 4001
$ws.Range("K17").Value = 4001
$ws.Range("M17").Value = -3829

# Row 38
$ws.Range("H38").Value = 23020.334
$ws.Range("J38").Value = 23020.334
$ws.Range("L38").Value = 23020.334
$ws.Range("N38").Value = -23966.334

# Row 43
$ws.Range("H43").Value = 17500
$ws.Range("J43").Value = 17500
$ws.Range("L43").Value = 17500
$ws.Range("N43").Value = -17798

# Row 49
$ws.Range("H49").Value = 5000
$ws.Range("I49").Value = 5000
$ws.Range("K49").Value = 5000
$ws.Range("M49").Value = -4770

# Row 113
$ws.Range("H113").Value = 1191.5454
$ws.Range("I113").Value = 973.375
$ws.Range("J113").Value = 1773.3334
$ws.Range("K113").Value = 2920.125
$ws.Range("L113").Value = 5320.0002
$ws.Range("M113").Value = -750.125
$ws.Range("N113").Value = -9660.0002

# Row 126
$ws.Range("H126").Value = 2626.6956
$ws.Range("J126").Value = 3089.6667
$ws.Range("L126").Value = 9269.000100000001
$ws.Range("N126").Value = -14209.0001

# Row 132
$ws.Range("H132").Value = 5321.65
$ws.Range("I132").Value = 7048.6294
$ws.Range("J132").Value = 1734.8462
$ws.Range("K132").Value = 21145.8882
$ws.Range("L132").Value = 5204.5386
$ws.Range("M132").Value = -18615.8882
$ws.Range("N132").Value = -10264.5386

